$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: add Wins / Losses / Ties labels, copying style from existing header cell (A1)
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Data rows 2-47: team win/loss/tie record (constant across all players on the roster)
for ($r = 2; $r -le 47; $r++) {
    $ws.Cells.Item($r, 30).Value = 69
    $ws.Cells.Item($r, 31).Value = 93
    $ws.Cells.Item($r, 32).Value = 0
}
